# Apply the edit described by the diff to the "GWP" worksheet:
#  - Select cell E11 (was A21)
#  - Replace static values in C10/D10/E10 with formulas

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GWP")

$ws.Range("C10").Formula = "=4.33+0.65"
$ws.Range("D10").Formula = "=D9+0.58"
$ws.Range("E10").Formula = "=E9+0.71"

$ws.Activate()
$ws.Range("E11").Select()
